$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.354.75"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.715.67"
$ws.Range("E3").Value = "  -1.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9968"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.30"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9975"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4870"
$ws.Range("E7").Value = "  -0.78%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2576"
$ws.Range("E8").Value = "  -3.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06165"
$ws.Range("E9").Value = "  -2.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.719.50"
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06945"
$ws.Range("E11").Value = "  -1.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.46"
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.475"
$ws.Range("E13").Value = "  -3.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5966"
$ws.Range("E14").Value = "  -2.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "76.40"
$ws.Range("E15").Value = "  -1.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9973"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.309.01"
$ws.Range("E17").Value = "  -0.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9970"
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("E19").Value = "  -4.61%  "
$ws.Range("E20").Value = "  -2.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.933.91"
$ws.Range("E21").Value = "  -1.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.408"
$ws.Range("E22").Value = "  -3.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.415"
$ws.Range("E23").Value = "  -3.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.046"
$ws.Range("E24").Value = "  -3.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.33"
$ws.Range("E25").Value = "  -3.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.15"
$ws.Range("E26").Value = "  -2.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.399"
$ws.Range("E27").Value = "  -1.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.725"
$ws.Range("E28").Value = "  -2.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "105.51"
$ws.Range("E29").Value = "  -2.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.868"
$ws.Range("E30").Value = "  -4.38%  "
$ws.Range("E31").Value = "  -1.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.611"
$ws.Range("E32").Value = "  -3.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04414"
$ws.Range("E33").Value = "  -3.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9964"
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.599"
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9920"
$ws.Range("E36").Value = "  -1.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6165"
$ws.Range("E37").Value = "  -3.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9346"
$ws.Range("E38").Value = "  +4.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.984"
$ws.Range("E39").Value = "  -1.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.372"
$ws.Range("E40").Value = "  -1.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9965"
$ws.Range("E41").Value = "  -0.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01471"
$ws.Range("E42").Value = "  -2.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.65"
$ws.Range("E43").Value = "  -3.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.430"
$ws.Range("E44").Value = "  +0.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3805"
$ws.Range("E45").Value = "  -2.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.812"
$ws.Range("E46").Value = "  -1.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1149"
$ws.Range("E47").Value = "  -3.25%  "
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("E49").Value = "  -0.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.692"
$ws.Range("E50").Value = "  -1.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.08"
$ws.Range("E51").Value = "  -1.44%  "
